# Add two new completed books to the "Completed" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# Fill the new string-valued cells in the same order the original author
# typed them in, so shared-string indices line up with the source file:
# A15, B15, E15, A16, B16, E16, G16, G15 (dates/F-column reuse existing
# shared strings so their order does not matter).
$ws.Cells.Item(15, 1).Value = "The Sports Gene"
$ws.Cells.Item(15, 2).Value = "David Epstein"
$ws.Cells.Item(15, 5).Value = "sports;science;genetics;expertise"

$ws.Cells.Item(16, 1).Value = "The Leadership Gap"
$ws.Cells.Item(16, 2).Value = "Lolly Daskal"
$ws.Cells.Item(16, 5).Value = "leadership;success;self-improvement"

$ws.Cells.Item(16, 7).Value = "6 Hrs 18 Mins"
$ws.Cells.Item(15, 7).Value = "10 Hrs 23 Mins"

# Dates (Start Date / Finish Date columns), formatted like the existing rows
# (built-in date format "m/d/yy" -> numFmtId 14, same style index as above).
$ws.Cells.Item(15, 3).NumberFormat = "m/d/yy"
$ws.Cells.Item(15, 3).Value = "1/20/2020"
$ws.Cells.Item(15, 4).NumberFormat = "m/d/yy"
$ws.Cells.Item(15, 4).Value = "1/23/2020"

$ws.Cells.Item(16, 3).NumberFormat = "m/d/yy"
$ws.Cells.Item(16, 3).Value = "1/23/2020"
$ws.Cells.Item(16, 4).NumberFormat = "m/d/yy"
$ws.Cells.Item(16, 4).Value = "1/24/2020"

# Type column (reuses the existing "Audio" shared string).
$ws.Cells.Item(15, 6).Value = "Audio"
$ws.Cells.Item(16, 6).Value = "Audio"

$ws.Range("G16").Select()
